$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.508.37'
$ws.Range("E2").Value = '  +1.42%  '

$ws.Range("D3").Value = '2.251.61'
$ws.Range("E3").Value = '  +1.03%  '

$ws.Range("E4").Value = '  +0.07%  '

$orig = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.78'
$ws.Range("D5").Style = $orig
$ws.Range("E5").Value = '  +2.12%  '

$orig = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.28'
$ws.Range("D6").Style = $orig
$ws.Range("E6").Value = '  +1.89%  '

$ws.Range("E7").Value = '  +1.42%  '

$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("E9").Value = '  +2.60%  '

$orig = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.26'
$ws.Range("D10").Style = $orig
$ws.Range("E10").Value = '  +3.72%  '

$ws.Range("E11").Value = '  +1.42%  '

$ws.Range("E12").Value = '  +3.44%  '

$orig = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.105'
$ws.Range("D13").Style = $orig
$ws.Range("E13").Value = '  +1.53%  '

$ws.Range("D14").Value = '2.311.20'
$ws.Range("E14").Value = '  +2.24%  '

$orig = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.840'
$ws.Range("D15").Style = $orig
$ws.Range("E15").Value = '  +4.32%  '

$ws.Range("E16").Value = '  +2.95%  '

$ws.Range("D17").Value = '44.232.07'
$ws.Range("E17").Value = '  +1.20%  '

$ws.Range("D18").Value = '0.0₃0966'
$ws.Range("E18").Value = '  +2.22%  '

$ws.Range("E19").Value = '  +5.82%  '

$orig = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.21'
$ws.Range("D20").Style = $orig
$ws.Range("E20").Value = '  +2.79%  '

$orig = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '66.08'
$ws.Range("D21").Style = $orig
$ws.Range("E21").Value = '  +2.91%  '

$orig = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.04'
$ws.Range("D22").Style = $orig
$ws.Range("E22").Value = '  +1.25%  '

$orig = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.01'
$ws.Range("D23").Style = $orig
$ws.Range("E23").Value = '  +4.38%  '

$ws.Range("E24").Value = '  +5.16%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("E26").Value = '  +5.97%  '

$ws.Range("E27").Value = '  +1.96%  '

$orig = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '38.02'
$ws.Range("D28").Style = $orig
$ws.Range("E28").Value = '  +6.56%  '

$orig = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.00'
$ws.Range("D29").Style = $orig
$ws.Range("E29").Value = '  +2.92%  '

$orig = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.13'
$ws.Range("D30").Style = $orig
$ws.Range("E30").Value = '  +2.14%  '

$orig = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.45'
$ws.Range("D31").Style = $orig

$orig = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0802'
$ws.Range("D32").Style = $orig
$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("E33").Value = '  +0.86%  '

$ws.Range("E34").Value = '  -1.16%  '

$orig = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.120'
$ws.Range("D35").Style = $orig
$ws.Range("E35").Value = '  +3.12%  '

$ws.Range("E36").Value = '  +3.29%  '

$ws.Range("E37").Value = '  +3.41%  '

$orig = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.46'
$ws.Range("D38").Style = $orig
$ws.Range("E38").Value = '  +6.34%  '

$orig = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.59'
$ws.Range("D39").Style = $orig
$ws.Range("E39").Value = '  +1.05%  '

$ws.Range("E40").Value = '  +1.50%  '

$ws.Range("E41").Value = '  +3.50%  '

$ws.Range("E42").Value = '  +0.15%  '

$ws.Range("D43").Value = '1.755.19'
$ws.Range("E43").Value = '  +1.13%  '

$ws.Range("E44").Value = '  +6.03%  '

$orig = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '81.02'
$ws.Range("D45").Style = $orig
$ws.Range("E45").Value = '  -2.23%  '

$orig = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.18'
$ws.Range("D46").Style = $orig
$ws.Range("E46").Value = '  +1.93%  '

$orig = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '71.07'
$ws.Range("D47").Style = $orig
$ws.Range("E47").Value = '  +5.29%  '

$orig = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.69'
$ws.Range("D48").Style = $orig
$ws.Range("E48").Value = '  +5.27%  '

$ws.Range("E49").Value = '  +3.34%  '

$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$orig = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.87'
$ws.Range("D50").Style = $orig
$ws.Range("E50").Value = '  +0.36%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$orig = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.59'
$ws.Range("D51").Style = $orig
$ws.Range("E51").Value = '  +6.64%  '
